$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(68,3).Value = "战报不存在"
$ws.Cells.Item(68,1).Value = "reportNotExist"
$ws.Cells.Item(68,2).Value = 567
$ws.Rows.Item(68).RowHeight = 20

$ws.Cells.Item(69,3).Value = "龙未处于空闲状态"
$ws.Cells.Item(69,1).Value = "dragonIsNotFree"
$ws.Cells.Item(69,2).Value = 568
$ws.Rows.Item(69).RowHeight = 20

$ws.Cells.Item(70,3).Value = "所选择的龙已经阵亡"
$ws.Cells.Item(70,1).Value = "dragonSelectedIsDead"
$ws.Cells.Item(70,2).Value = 569
$ws.Rows.Item(70).RowHeight = 20

$ws.Cells.Item(71,3).Value = "没有龙驻防在城墙"
$ws.Cells.Item(71,1).Value = "noDragonInDefenceStatus"
$ws.Cells.Item(71,2).Value = 570
$ws.Rows.Item(71).RowHeight = 20

$ws.Cells.Item(72,3).Value = "没有足够的出售队列"
$ws.Cells.Item(72,1).Value = "sellQueueNotEnough"
$ws.Cells.Item(72,2).Value = 571
$ws.Rows.Item(72).RowHeight = 20

$ws.Cells.Item(73,3).Value = "玩家资源不足"
$ws.Cells.Item(73,1).Value = "resourceNotEnough"
$ws.Cells.Item(73,2).Value = 572
$ws.Rows.Item(73).RowHeight = 20

$ws.Cells.Item(74,3).Value = "马车数量不足"
$ws.Cells.Item(74,1).Value = "cartNotEnough"
$ws.Cells.Item(74,2).Value = 573
$ws.Rows.Item(74).RowHeight = 20

$ws.Cells.Item(75,3).Value = "商品不存在"
$ws.Cells.Item(75,1).Value = "sellItemNotExist"
$ws.Cells.Item(75,2).Value = 574
$ws.Rows.Item(75).RowHeight = 20

$ws.Cells.Item(76,3).Value = "银币不足"
$ws.Cells.Item(76,1).Value = "coinNotEnough"
$ws.Cells.Item(76,2).Value = 575
$ws.Rows.Item(76).RowHeight = 20

$ws.Cells.Item(77,3).Value = "商品还未卖出"
$ws.Cells.Item(77,1).Value = "sellItemNotSold"
$ws.Cells.Item(77,2).Value = 576
$ws.Rows.Item(77).RowHeight = 20

$ws.Cells.Item(78,3).Value = "您未出售此商品"
$ws.Cells.Item(78,1).Value = "sellItemNotBelongsToYou"
$ws.Cells.Item(78,2).Value = 577
$ws.Rows.Item(78).RowHeight = 20

$ws.Cells.Item(79,3).Value = "商品已经售出"
$ws.Cells.Item(79,1).Value = "sellItemAlreadySold"
$ws.Cells.Item(79,2).Value = 578
$ws.Rows.Item(79).RowHeight = 20

$ws.Range("C79").Select()
